# "atualizacao 2a RQ 2024-2028" - update market growth-rate assumptions
# and extend the projection horizon from 2050 through 2060.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Years 2024-2050 (rows 16-42): refreshed growth-rate assumptions.
# 2024-2027 get distinct new rates; 2028 onward flatten to 2.6%.
$rates = @{
    16 = 0.022  # 2024
    17 = 0.022  # 2025
    18 = 0.023  # 2026
    19 = 0.025  # 2027
}
foreach ($row in $rates.Keys) {
    $ws.Range("B$row").Value = $rates[$row]
}

for ($r = 20; $r -le 42; $r++) {
    $ws.Range("B$r").Value = 0.026
}

# Extend the series with new rows for years 2051-2060 (rows 43-52),
# continuing the flat 2.6% growth-rate assumption.
$year = 2051
for ($r = 43; $r -le 52; $r++) {
    $ws.Range("A$r").Value = $year
    $ws.Range("B$r").Value = 0.026
    $year++
}

# Reflect the cursor/selection position left in the saved view.
$ws.Range("B20").Select()

$wb.Save()
